$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Most of these are on - Opengameart.org" -> " Opengameart.org"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Most of these are on - Opengameart.org", $true, $false, $false, $false, $false, $true, 1, $false, " Opengameart.org", 2) | Out-Null

# ------------------------------------------------------------------
# 2. "Minotaur - PixElthen" -> "Minotaur Sprite - PixElthen"
#    (insert a new "Sprite " word between "Minotaur " and "- PixElthen")
# ------------------------------------------------------------------
$fr = $d.Content
$fr.Find.Execute("Minotaur ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins = $d.Range($fr.End, $fr.End)
$ins.InsertAfter("Sprite ")

# ------------------------------------------------------------------
# 3. "Grey Minotaur - Merry Dream Games" -> "Grey Minotaur Sprite - Merry Dream Games"
#    (insert " Sprite" right after "Grey Minotaur")
# ------------------------------------------------------------------
$fr2 = $d.Content
$fr2.Find.Execute("Grey Minotaur", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins2 = $d.Range($fr2.End, $fr2.End)
$ins2.InsertAfter(" Sprite")

# ------------------------------------------------------------------
# 4. "Cave background - Spring" -> "Cave background <en dash> Spring"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Cave background - Spring", $true, $false, $false, $false, $false, $true, 1, $false, "Cave background " + [char]0x2013 + " Spring", 2) | Out-Null

# ------------------------------------------------------------------
# 5. Relocate the "_GoBack" bookmark: delete it here, it will be
#    re-added at its new home once the trailing paragraphs exist.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 6. Add two new paragraphs after "Cave background ... Spring":
#      "Coin Sprite - irmirx "
#      "Title, start and end game buttons - Joseph"
# ------------------------------------------------------------------
$caveBackgroundPara = $d.Paragraphs($d.Paragraphs.Count)
$caveBackgroundPara.Range.InsertParagraphAfter()

$coinPara = $d.Paragraphs($d.Paragraphs.Count)
$coinPara.Range.InsertBefore("Coin Sprite " + [char]0x2013 + " irmirx ")

$coinPara2 = $d.Paragraphs($d.Paragraphs.Count)
$coinPara2.Range.InsertParagraphAfter()

$titlePara = $d.Paragraphs($d.Paragraphs.Count)
$titleStart = "Title, start and end game buttons"
$titleTail = " " + [char]0x2013 + " Joseph"
# Insert the full final text (including the tail) first -- inserting the
# bookmark exactly at the end of the document's content has proven to
# misplace it, so we make sure there is trailing text past the bookmark
# point before we add it.
$titlePara.Range.InsertBefore($titleStart + $titleTail)

# ------------------------------------------------------------------
# 7. Re-seat the "_GoBack" bookmark between "...buttons" and " - Joseph"
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs($d.Paragraphs.Count)
$bmPos = $titlePara.Range.Start + $titleStart.Length
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))
